$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 (Arkansas) updates ---
$ws.Range("C9").Value = 24253
$ws.Range("D9").Value = 292
$ws.Range("E9").Value = 5239
$ws.Range("G9").Value = 25.23
$ws.Range("H9").Value = 25.18
$ws.Range("K9").Value = 20762
$ws.Range("L9").Value = 278

# --- Row 10 (California - San Diego) updates ---
$ws.Range("B10").Value = 44018
$ws.Range("C10").Value = 17000
$ws.Range("E10").Value = 615
$ws.Range("G10").Value = 4.55
$ws.Range("K10").Value = 13514

# --- Row 31 (Washington) updates ---
$ws.Range("B31").Value = 44018
$ws.Range("C31").Value = 36985
$ws.Range("D31").Value = 1370
$ws.Range("E31").Value = 1455
$ws.Range("K31").Value = 26515

# --- New row 42 (Iowa) ---
$ws.Range("A42").Value = "Iowa"
$ws.Range("I42").Value = $false
$ws.Range("J42").Value = $false
$ws.Range("M42").Value = 109911
$ws.Range("N42").Value = 3.51

$quote = [char]39
$errMsg = "An error occurred. ... ValueError(" + $quote + "Unable to parse " + [char]34 + "Reported Deaths In Adair : No Data" + [char]34 + " as int" + $quote + ")"
$ws.Range("O42").Value = $errMsg
